$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header date bumped by one day (Excel serial 45310 -> 45311)
$ws.Range("A1").Value = 45311

# Price list (column D) updated for rows 33-38
$ws.Range("D33").Value = 94
$ws.Range("D34").Value = 75
$ws.Range("D35").Value = 70
$ws.Range("D36").Value = 161
$ws.Range("D37").Value = 115
$ws.Range("D38").Value = 101
